$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.950.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.552.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.66%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.550.95"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.159.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("E14").Value = "  +4.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.564.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.865.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.577"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.699.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.72%  "
$ws.Range("E25").Value = "  +3.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.562.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +22.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.79%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.64%  "
$ws.Range("E39").Value = "  +8.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.827"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.82%  "
$ws.Range("E47").Value = "  +11.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.468.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.61%  "
